$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value. All values are written as
# literal Text (matching the source inlineStr cells) regardless of whether
# the string happens to look like a number.
$updates = @(
    @{Cell="D2"; Value="26.483.33"},
    @{Cell="E2"; Value="  -7.74%  "},
    @{Cell="D3"; Value="1.678.43"},
    @{Cell="E3"; Value="  -6.90%  "},
    @{Cell="D4"; Value="1.004"},
    @{Cell="E4"; Value="  +0.18%  "},
    @{Cell="D5"; Value="216.29"},
    @{Cell="E5"; Value="  -6.56%  "},
    @{Cell="D6"; Value="1.005"},
    @{Cell="E6"; Value="  +0.16%  "},
    @{Cell="D7"; Value="0.4962"},
    @{Cell="E7"; Value="  -16.61%  "},
    @{Cell="D8"; Value="0.2593"},
    @{Cell="E8"; Value="  -6.71%  "},
    @{Cell="D9"; Value="21.71"},
    @{Cell="E9"; Value="  -6.79%  "},
    @{Cell="D10"; Value="0.06159"},
    @{Cell="E10"; Value="  -9.88%  "},
    @{Cell="D11"; Value="0.07276"},
    @{Cell="E11"; Value="  -3.51%  "},
    @{Cell="D12"; Value="1.663.84"},
    @{Cell="E12"; Value="  -5.27%  "},
    @{Cell="D13"; Value="4.421"},
    @{Cell="E13"; Value="  -7.23%  "},
    @{Cell="D14"; Value="0.5723"},
    @{Cell="E14"; Value="  -8.14%  "},
    @{Cell="D15"; Value="1.907.22"},
    @{Cell="E15"; Value="  -6.86%  "},
    @{Cell="E16"; Value="  -12.71%  "},
    @{Cell="D17"; Value="64.15"},
    @{Cell="E17"; Value="  -14.99%  "},
    @{Cell="D18"; Value="26.485.47"},
    @{Cell="E18"; Value="  -7.57%  "},
    @{Cell="D19"; Value="4.973"},
    @{Cell="E19"; Value="  -9.27%  "},
    @{Cell="D20"; Value="1.005"},
    @{Cell="E20"; Value="  +0.17%  "},
    @{Cell="E21"; Value="  -6.25%  "},
    @{Cell="D22"; Value="183.30"},
    @{Cell="E22"; Value="  -12.63%  "},
    @{Cell="D23"; Value="6.155"},
    @{Cell="E23"; Value="  -10.28%  "},
    @{Cell="D24"; Value="1.005"},
    @{Cell="E24"; Value="  +0.21%  "},
    @{Cell="D25"; Value="144.25"},
    @{Cell="E25"; Value="  -6.47%  "},
    @{Cell="D26"; Value="7.461"},
    @{Cell="E26"; Value="  -4.99%  "},
    @{Cell="D27"; Value="0.1123"},
    @{Cell="E27"; Value="  -11.67%  "},
    @{Cell="D28"; Value="15.41"},
    @{Cell="E28"; Value="  -5.90%  "},
    @{Cell="E29"; Value="  -8.81%  "},
    @{Cell="D30"; Value="0.05685"},
    @{Cell="E30"; Value="  -8.23%  "},
    @{Cell="D31"; Value="1.318"},
    @{Cell="E31"; Value="  -7.20%  "},
    @{Cell="D32"; Value="3.465"},
    @{Cell="E32"; Value="  -8.28%  "},
    @{Cell="D33"; Value="3.449"},
    @{Cell="E33"; Value="  -7.86%  "},
    @{Cell="D34"; Value="1.625"},
    @{Cell="E34"; Value="  -5.44%  "},
    @{Cell="D35"; Value="1.002"},
    @{Cell="E35"; Value="  -5.88%  "},
    @{Cell="D36"; Value="2.366"},
    @{Cell="E36"; Value="  -5.07%  "},
    @{Cell="D37"; Value="0.5879"},
    @{Cell="E37"; Value="  -8.00%  "},
    @{Cell="D38"; Value="2.631"},
    @{Cell="E38"; Value="  -3.00%  "},
    @{Cell="D39"; Value="0.01583"},
    @{Cell="E39"; Value="  -7.51%  "},
    @{Cell="D40"; Value="1.067.40"},
    @{Cell="E40"; Value="  -5.65%  "},
    @{Cell="D41"; Value="5.863"},
    @{Cell="E41"; Value="  -9.13%  "},
    @{Cell="D42"; Value="0.8487"},
    @{Cell="E42"; Value="  -2.66%  "},
    @{Cell="D43"; Value="1.003"},
    @{Cell="E43"; Value="  -0.08%  "},
    @{Cell="D44"; Value="97.98"},
    @{Cell="E44"; Value="  -2.67%  "},
    @{Cell="D45"; Value="1.835.73"},
    @{Cell="E45"; Value="  -6.33%  "},
    @{Cell="D46"; Value="55.98"},
    @{Cell="E46"; Value="  -7.46%  "},
    @{Cell="B47"; Value="BabyDogeCoin"},
    @{Cell="C47"; Value="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"},
    @{Cell="D47"; Value="0.00000000105"},
    @{Cell="E47"; Value="  -6.34%  "},
    @{Cell="B48"; Value="Frax"},
    @{Cell="C48"; Value="https://coinranking.com/coin/KfWtaeV1W+frax-frax"},
    @{Cell="D48"; Value="1.005"},
    @{Cell="E48"; Value="  -0.14%  "},
    @{Cell="D49"; Value="8.042"},
    @{Cell="E49"; Value="  -3.35%  "},
    @{Cell="D50"; Value="0.4309"},
    @{Cell="E50"; Value="  -3.95%  "},
    @{Cell="D51"; Value="0.05183"},
    @{Cell="E51"; Value="  -5.15%  "}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $u.Value
    $rng.Style = "Normal"
}